# Auto-generated edit script applying numeric corrections to H:N columns
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# matching the target diff for Sophia_Profits workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 173
$ws.Cells.Item(5, 10).Value = 240
$ws.Cells.Item(5, 12).Value = 240
$ws.Cells.Item(5, 14).Value = -470
$ws.Cells.Item(62, 8).Value = 1948
$ws.Cells.Item(62, 9).Value = 1899.5
$ws.Cells.Item(62, 10).Value = 1996.5
$ws.Cells.Item(62, 11).Value = 1899.5
$ws.Cells.Item(62, 12).Value = 1996.5
$ws.Cells.Item(62, 13).Value = -1275.5
$ws.Cells.Item(62, 14).Value = -3244.5
$ws.Cells.Item(65, 8).Value = 1948
$ws.Cells.Item(65, 9).Value = 1899.5
$ws.Cells.Item(65, 10).Value = 1996.5
$ws.Cells.Item(65, 11).Value = 9497.5
$ws.Cells.Item(65, 12).Value = 9982.5
$ws.Cells.Item(65, 13).Value = -6377.5
$ws.Cells.Item(65, 14).Value = -16222.5
$ws.Cells.Item(86, 8).Value = 8622.5
$ws.Cells.Item(86, 9).Value = 8500
$ws.Cells.Item(86, 11).Value = 8500
$ws.Cells.Item(86, 13).Value = -7377
$ws.Cells.Item(89, 8).Value = 8622.5
$ws.Cells.Item(89, 9).Value = 8500
$ws.Cells.Item(89, 11).Value = 42500
$ws.Cells.Item(89, 13).Value = -36884
$ws.Cells.Item(92, 8).Value = 510.15384
$ws.Cells.Item(92, 9).Value = 597.8182
$ws.Cells.Item(92, 10).Value = 28
$ws.Cells.Item(92, 11).Value = 597.8182
$ws.Cells.Item(92, 12).Value = 28
$ws.Cells.Item(92, 13).Value = 650.1818
$ws.Cells.Item(92, 14).Value = -2524
$ws.Cells.Item(125, 8).Value = 932.6667
$ws.Cells.Item(125, 9).Value = 899
$ws.Cells.Item(125, 10).Value = 1000
$ws.Cells.Item(125, 11).Value = 8091
$ws.Cells.Item(125, 12).Value = 9000
$ws.Cells.Item(125, 13).Value = -5631
$ws.Cells.Item(125, 14).Value = -13920
$ws.Cells.Item(132, 8).Value = 9365.166999999999
$ws.Cells.Item(132, 9).Value = 9238.4
$ws.Cells.Item(132, 11).Value = 27715.2
$ws.Cells.Item(132, 13).Value = -25185.2
$ws.Cells.Item(137, 8).Value = 3521.6667
$ws.Cells.Item(137, 9).Value = 3126
$ws.Cells.Item(137, 10).Value = 5500
$ws.Cells.Item(137, 11).Value = 9378
$ws.Cells.Item(137, 12).Value = 16500
$ws.Cells.Item(137, 13).Value = -6828
$ws.Cells.Item(137, 14).Value = -21600
$ws.Cells.Item(141, 8).Value = 950
$ws.Cells.Item(141, 9).Value = 950
$ws.Cells.Item(141, 11).Value = 2850
$ws.Cells.Item(141, 13).Value = 2330

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6740.231
$ws.Cells.Item(32, 9).Value = 6740.231
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 6740.231
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -6453.231
$ws.Cells.Item(32, 14).Value = $null
$ws.Cells.Item(61, 8).Value = 2942
$ws.Cells.Item(61, 9).Value = 2942
$ws.Cells.Item(61, 11).Value = 2942
$ws.Cells.Item(61, 13).Value = -2730
$ws.Cells.Item(97, 8).Value = 846.6667
$ws.Cells.Item(97, 9).Value = 846.6667
$ws.Cells.Item(97, 11).Value = 846.6667
$ws.Cells.Item(97, 13).Value = -350.6667
$ws.Cells.Item(122, 8).Value = 1417.1666
$ws.Cells.Item(122, 9).Value = 1137.8
$ws.Cells.Item(122, 10).Value = 2814
$ws.Cells.Item(122, 11).Value = 3413.4
$ws.Cells.Item(122, 12).Value = 8442
$ws.Cells.Item(122, 13).Value = -963.3999999999996
$ws.Cells.Item(122, 14).Value = -13342
$ws.Cells.Item(132, 8).Value = 3771
$ws.Cells.Item(132, 9).Value = 3771
$ws.Cells.Item(132, 11).Value = 11313
$ws.Cells.Item(132, 13).Value = -8783
$ws.Cells.Item(136, 8).Value = 2942
$ws.Cells.Item(136, 9).Value = 2942
$ws.Cells.Item(136, 11).Value = 8826
$ws.Cells.Item(136, 13).Value = -6276

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 460.5
$ws.Cells.Item(80, 9).Value = 53
$ws.Cells.Item(80, 10).Value = 705
$ws.Cells.Item(80, 11).Value = 53
$ws.Cells.Item(80, 12).Value = 705
$ws.Cells.Item(80, 13).Value = 945
$ws.Cells.Item(80, 14).Value = -2701
$ws.Cells.Item(83, 8).Value = 460.5
$ws.Cells.Item(83, 9).Value = 53
$ws.Cells.Item(83, 10).Value = 705
$ws.Cells.Item(83, 11).Value = 265
$ws.Cells.Item(83, 12).Value = 3525
$ws.Cells.Item(83, 13).Value = 4727
$ws.Cells.Item(83, 14).Value = -13509
$ws.Cells.Item(134, 8).Value = 3698.2
$ws.Cells.Item(134, 9).Value = 3698.2
$ws.Cells.Item(134, 11).Value = 11094.6
$ws.Cells.Item(134, 13).Value = -8559.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 183
$ws.Cells.Item(2, 9).Value = 202
$ws.Cells.Item(2, 10).Value = 145
$ws.Cells.Item(2, 11).Value = 202
$ws.Cells.Item(2, 12).Value = 145
$ws.Cells.Item(2, 13).Value = -89
$ws.Cells.Item(2, 14).Value = -371
$ws.Cells.Item(16, 8).Value = 1000000
$ws.Cells.Item(16, 9).Value = 1000000
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1000000
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -999713
$ws.Cells.Item(16, 14).Value = $null
$ws.Cells.Item(47, 8).Value = 5000
$ws.Cells.Item(47, 10).Value = 5000
$ws.Cells.Item(47, 12).Value = 5000
$ws.Cells.Item(47, 14).Value = -6132
$ws.Cells.Item(74, 8).Value = 59984.668
$ws.Cells.Item(74, 10).Value = 59984.668
$ws.Cells.Item(74, 12).Value = 59984.668
$ws.Cells.Item(74, 14).Value = -61732.668
$ws.Cells.Item(77, 8).Value = 59984.668
$ws.Cells.Item(77, 10).Value = 59984.668
$ws.Cells.Item(77, 12).Value = 179954.004
$ws.Cells.Item(77, 14).Value = -188690.004
$ws.Cells.Item(113, 8).Value = 1000000
$ws.Cells.Item(113, 9).Value = 1000000
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 1000000
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -997830
$ws.Cells.Item(113, 14).Value = $null
$ws.Cells.Item(132, 8).Value = 146542
$ws.Cells.Item(132, 9).Value = 203159
$ws.Cells.Item(132, 11).Value = 609477
$ws.Cells.Item(132, 13).Value = -606947
$ws.Cells.Item(134, 8).Value = 3000
$ws.Cells.Item(134, 9).Value = 3000
$ws.Cells.Item(134, 11).Value = 9000
$ws.Cells.Item(134, 13).Value = -6465

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(102, 8).Value = 5333
$ws.Cells.Item(102, 10).Value = 5333
$ws.Cells.Item(102, 12).Value = 15999
$ws.Cells.Item(102, 14).Value = -20867

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3453.0908
$ws.Cells.Item(126, 9).Value = 3163.3333
$ws.Cells.Item(126, 11).Value = 9489.999899999999
$ws.Cells.Item(126, 13).Value = -7019.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 7635.647
$ws.Cells.Item(22, 9).Value = 9104.5
$ws.Cells.Item(22, 10).Value = 6330
$ws.Cells.Item(22, 11).Value = 9104.5
$ws.Cells.Item(22, 12).Value = 6330
$ws.Cells.Item(22, 13).Value = -8809.5
$ws.Cells.Item(22, 14).Value = -6920
$ws.Cells.Item(27, 8).Value = 7635.647
$ws.Cells.Item(27, 9).Value = 9104.5
$ws.Cells.Item(27, 10).Value = 6330
$ws.Cells.Item(27, 11).Value = 9104.5
$ws.Cells.Item(27, 12).Value = 6330
$ws.Cells.Item(27, 13).Value = -8997.5
$ws.Cells.Item(27, 14).Value = -6544
$ws.Cells.Item(55, 8).Value = 887.5
$ws.Cells.Item(55, 9).Value = 887.5
$ws.Cells.Item(55, 11).Value = 887.5
$ws.Cells.Item(55, 13).Value = -714.5
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).Value = $null
$ws.Cells.Item(132, 8).Value = 2534.3
$ws.Cells.Item(132, 9).Value = 1793.125
$ws.Cells.Item(132, 10).Value = 5499
$ws.Cells.Item(132, 11).Value = 5379.375
$ws.Cells.Item(132, 12).Value = 16497
$ws.Cells.Item(132, 13).Value = -2849.375
$ws.Cells.Item(132, 14).Value = -21557

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 89995
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 13).Value = $null
$ws.Cells.Item(67, 8).Value = 89995
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 13).Value = $null
$ws.Cells.Item(140, 8).Value = 75494
$ws.Cells.Item(140, 10).Value = 75494
$ws.Cells.Item(140, 12).Value = 75494
$ws.Cells.Item(140, 14).Value = -85854
